$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CapitalRemittance")

# Add the new "Folio No*" column header in H1
$ws.Range("H1").Value = "Folio No*"

# Update the selection to reflect the new active cell (H2) as in the saved file
$ws.Range("H2").Select()
